# Compromise_Assessment.xlsx - refactor validate.py, refactor make.py
# Remove the Linux checklist section (old rows 15-28) and re-curate the
# Windows checklist items (rows 2-15), then normalize borders across the
# whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Drop the entire "Linux" block (old rows 16-28); old row 15 (the first
#    Linux data row) gets repurposed below to hold the final Windows item.
$ws.Range("A16:A28").EntireRow.Delete()

# 2) Re-point row 15 at the Windows checklist (it still carries the old
#    "Linux" header values/styles inherited from the former row 15).
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Windows"
$ws.Cells.Item(15, 3).Value = "Check Malicious Files"

# 3) Re-curate the checklist item names in column C (rows 4-15 shuffle /
#    change relative to the pre-edit sheet).
$ws.Cells.Item(4, 3).Value  = "Check the service Everyone Permission"
$ws.Cells.Item(5, 3).Value  = "Suspicious Directory"
$ws.Cells.Item(6, 3).Value  = "Visual Basic for Applications"
$ws.Cells.Item(7, 3).Value  = "Startup files"
$ws.Cells.Item(8, 3).Value  = "Living off the Land"
$ws.Cells.Item(9, 3).Value  = "Event Fies Check"
$ws.Cells.Item(10, 3).Value = "schedule Task"
$ws.Cells.Item(11, 3).Value = "Defender and Realtime monitoring"
$ws.Cells.Item(12, 3).Value = "Third Party Application"
$ws.Cells.Item(13, 3).Value = "Auto Runs"
$ws.Cells.Item(14, 3).Value = "Powershell History"

# 4) Row 14 height shrinks to match the other body rows.
$ws.Rows.Item(14).RowHeight = 18.75

# 5) Apply a uniform thin border across the whole surviving table.
$table = $ws.Range("A1:C15")
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2

# B5 keeps the explicit-black font color left over from editing (matches
# the source diff).
$ws.Cells.Item(5, 2).Font.Color = 0

# 6) Sheet view: selection/scroll position moved while editing.
$ws.Range("H8").Select()
$excel.ActiveWindow.ScrollRow = 5
